$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.858.35'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.635.56'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.02'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5092'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2582'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06423'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07800'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.271'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.641.07'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.27'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7660'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.21'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.859.38'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.33'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.945'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.162'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.788'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '138.53'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1231'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.852'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.53'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.239'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04949'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.300'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.248'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.87%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.387'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9036'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.574'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.57%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.134.58'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.83%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9969'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.465'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8003'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.47'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4265'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.777'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05077'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.24%  '
